# Scheduled-runner style refresh of cached price/profit figures across the
# per-job profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
# Only numeric "price/profit" columns (H..N) are refreshed per affected row;
# a few cells had no valid figure for this run and are cleared entirely.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 34565
$ws.Range("I32").Value = 1700
$ws.Range("J32").Value = 50997.5
$ws.Range("K32").Value = 1700
$ws.Range("L32").Value = 50997.5
$ws.Range("M32").Value = -1374
$ws.Range("N32").Value = -51649.5

$ws.Range("H58").Value = 3506.3333
$ws.Range("I58").Value = 257.375
$ws.Range("K58").Value = 772.125
$ws.Range("M58").Value = -622.125

$ws.Range("H76").Value = 50004264
$ws.Range("J76").Value = 4768.25
$ws.Range("L76").Value = 4768.25
$ws.Range("N76").Value = -5398.25

$ws.Range("H79").Value = 50004264
$ws.Range("J79").Value = 4768.25
$ws.Range("L79").Value = 4768.25
$ws.Range("N79").Value = -6952.25

$ws.Range("H106").Value = 696543.25
$ws.Range("I106").Value = 927641.3
$ws.Range("K106").Value = 927641.3
$ws.Range("M106").Value = -927010.3

$ws.Range("H111").Value = 502.13635
$ws.Range("I111").Value = 424.10526
$ws.Range("K111").Value = 1272.31578
$ws.Range("M111").Value = 1794.68422

$ws.Range("H133").Value = 98568.42999999999
$ws.Range("J133").Value = 98568.42999999999
$ws.Range("L133").Value = 98568.42999999999
$ws.Range("N133").Value = -108688.43

$ws.Range("H134").Value = 82642.44500000001
$ws.Range("J134").Value = 82642.44500000001
$ws.Range("L134").Value = 82642.44500000001
$ws.Range("N134").Value = -92782.44500000001

$ws.Range("H136").Value = 71437.55499999999
$ws.Range("J136").Value = 71437.55499999999
$ws.Range("L136").Value = 71437.55499999999
$ws.Range("N136").Value = -81637.55499999999

$ws.Range("H139").Value = 69470.7
$ws.Range("J139").Value = 69470.7
$ws.Range("L139").Value = 69470.7
$ws.Range("N139").Value = -79750.7

$ws.Range("H140").Value = 91051.82000000001
$ws.Range("J140").Value = 91557.10000000001
$ws.Range("L140").Value = 91557.10000000001
$ws.Range("N140").Value = -101917.1

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H6").Value = 69625
$ws.Range("I6").Value = 69625
$ws.Range("J6").Value = 0
$ws.Range("K6").Value = 69625
$ws.Range("L6").Value = 0
$ws.Range("N6").Value = -69452
$ws.Range("M6").ClearContents()

$ws.Range("H32").Value = 5440.0884
$ws.Range("I32").Value = 2977.9614
$ws.Range("K32").Value = 2977.9614
$ws.Range("M32").Value = -2690.9614

$ws.Range("H43").Value = 67775
$ws.Range("J43").Value = 79995
$ws.Range("L43").Value = 79995
$ws.Range("N43").Value = -80621

$ws.Range("H46").Value = 26332
$ws.Range("J46").Value = 29498
$ws.Range("L46").Value = 29498
$ws.Range("N46").Value = -30136

$ws.Range("H74").Value = 5055.8887
$ws.Range("I74").Value = 2692.7693
$ws.Range("J74").Value = 11200
$ws.Range("K74").Value = 2692.7693
$ws.Range("L74").Value = 11200
$ws.Range("M74").Value = -1818.7693
$ws.Range("N74").Value = -12948

$ws.Range("H77").Value = 5055.8887
$ws.Range("I77").Value = 2692.7693
$ws.Range("J77").Value = 11200
$ws.Range("K77").Value = 13463.8465
$ws.Range("L77").Value = 56000
$ws.Range("M77").Value = -9095.8465
$ws.Range("N77").Value = -64736

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H132").Value = 24620.428
$ws.Range("J132").Value = 24620.428
$ws.Range("L132").Value = 24620.428
$ws.Range("N132").Value = -34740.428

$ws.Range("H134").Value = 5187.4546
$ws.Range("I134").Value = 3605.3572
$ws.Range("K134").Value = 10816.0716
$ws.Range("M134").Value = -8281.071599999999

$ws.Range("H135").Value = 91997.375
$ws.Range("J135").Value = 91997.375
$ws.Range("L135").Value = 91997.375
$ws.Range("N135").Value = -102137.375

$ws.Range("H138").Value = 86738
$ws.Range("J138").Value = 86738
$ws.Range("L138").Value = 86738
$ws.Range("N138").Value = -97018

$ws.Range("H140").Value = 82132.14
$ws.Range("J140").Value = 82132.14
$ws.Range("L140").Value = 82132.14
$ws.Range("N140").Value = -92492.14

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H6").Value = 251174.12
$ws.Range("I6").Value = 144199
$ws.Range("J6").Value = 1000000
$ws.Range("K6").Value = 144199
$ws.Range("L6").Value = 1000000
$ws.Range("M6").Value = -144086
$ws.Range("N6").Value = -1000226

$ws.Range("H50").Value = 15000
$ws.Range("J50").Value = 0
$ws.Range("L50").Value = 0
$ws.Range("N50").ClearContents()

$ws.Range("H62").Value = 2890
$ws.Range("I62").Value = 2890
$ws.Range("K62").Value = 2890
$ws.Range("M62").Value = -2266

$ws.Range("H65").Value = 2890
$ws.Range("I65").Value = 2890
$ws.Range("K65").Value = 14450
$ws.Range("M65").Value = -11330

$ws.Range("H138").Value = 83428.875
$ws.Range("J138").Value = 83428.875
$ws.Range("L138").Value = 83428.875
$ws.Range("N138").Value = -93708.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H10").Value = 194.21053
$ws.Range("J10").Value = 164.75
$ws.Range("L10").Value = 494.25
$ws.Range("N10").Value = -772.25

$ws.Range("H18").Value = 15064.286
$ws.Range("I18").Value = 20290
$ws.Range("K18").Value = 60870
$ws.Range("M18").Value = -60701

$ws.Range("H107").Value = 1045.3704
$ws.Range("I107").Value = 1058.7273
$ws.Range("J107").Value = 1036.1875
$ws.Range("K107").Value = 3176.1819
$ws.Range("L107").Value = 3108.5625
$ws.Range("M107").Value = -1256.1819
$ws.Range("N107").Value = -6948.5625

$ws.Range("H121").Value = 715601.8
$ws.Range("J121").Value = 2501132.8
$ws.Range("L121").Value = 7503398.399999999
$ws.Range("N121").Value = -7506018.399999999

$ws.Range("H131").Value = 1144.2354
$ws.Range("I131").Value = 856.9167
$ws.Range("J131").Value = 1833.8
$ws.Range("K131").Value = 2570.7501
$ws.Range("L131").Value = 5501.4
$ws.Range("M131").Value = 2469.2499
$ws.Range("N131").Value = -15581.4

$ws.Range("H133").Value = 5000
$ws.Range("I133").Value = 5000
$ws.Range("J133").Value = 0
$ws.Range("K133").Value = 15000
$ws.Range("L133").Value = 0
$ws.Range("N133").Value = -9940
$ws.Range("M133").ClearContents()

$ws.Range("H134").Value = 92949.17999999999
$ws.Range("I134").Value = 92949.17999999999
$ws.Range("K134").Value = 278847.54
$ws.Range("M134").Value = -273777.54

$ws.Range("H139").Value = 1259.7142
$ws.Range("I139").Value = 1259.7142
$ws.Range("K139").Value = 3779.1426
$ws.Range("M139").Value = 1360.8574

$ws.Range("H140").Value = 2100
$ws.Range("I140").Value = 2100
$ws.Range("K140").Value = 6300
$ws.Range("M140").Value = -1120

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 11038.5
$ws.Range("J93").Value = 11038.5
$ws.Range("L93").Value = 11038.5
$ws.Range("N93").Value = -14782.5

$ws.Range("H135").Value = 99757.60000000001
$ws.Range("J135").Value = 99757.60000000001
$ws.Range("L135").Value = 99757.60000000001
$ws.Range("N135").Value = -109897.6

$ws.Range("H140").Value = 89547.37
$ws.Range("J140").Value = 89547.37
$ws.Range("L140").Value = 89547.37
$ws.Range("N140").Value = -99907.37

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 8125.875
$ws.Range("I19").Value = 7349.5
$ws.Range("J19").Value = 8902.25
$ws.Range("K19").Value = 7349.5
$ws.Range("L19").Value = 8902.25
$ws.Range("M19").Value = -7179.5
$ws.Range("N19").Value = -9242.25

$ws.Range("H40").Value = 3539001.5
$ws.Range("I40").Value = 81267.234
$ws.Range("J40").Value = 18522516
$ws.Range("K40").Value = 81267.234
$ws.Range("L40").Value = 18522516
$ws.Range("M40").Value = -81131.234
$ws.Range("N40").Value = -18522788

$ws.Range("H106").Value = 28624.285
$ws.Range("J106").Value = 28624.285
$ws.Range("L106").Value = 28624.285
$ws.Range("N106").Value = -31148.285

$ws.Range("H122").Value = 71575944
$ws.Range("I122").Value = 77081170
$ws.Range("J122").Value = 8005
$ws.Range("K122").Value = 231243510
$ws.Range("L122").Value = 24015
$ws.Range("M122").Value = -231241060
$ws.Range("N122").Value = -28915

$ws.Range("H136").Value = 4167.4165
$ws.Range("I136").Value = 4028.842
$ws.Range("J136").Value = 4694
$ws.Range("K136").Value = 12086.526
$ws.Range("L136").Value = 14082
$ws.Range("M136").Value = -9536.526
$ws.Range("N136").Value = -19182

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H33").Value = 25000
$ws.Range("J33").Value = 0
$ws.Range("L33").Value = 0
$ws.Range("N33").ClearContents()

$ws.Range("H36").Value = 25000
$ws.Range("J36").Value = 0
$ws.Range("L36").Value = 0
$ws.Range("N36").ClearContents()

$ws.Range("H136").Value = 2333.111
$ws.Range("I136").Value = 1399.6
$ws.Range("K136").Value = 4198.799999999999
$ws.Range("M136").Value = -1648.799999999999

